$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# The category list is re-ordered: the "Trays" category (formerly the last
# row, A19) moves up to row 9 - the slot vacated by removing "Fresh Juices".
# Every category that used to sit between them (Cocktails ... Shisha) shifts
# down by one row, the old "Fresh Juices" entry is relabelled to the new
# "Fresh Juice" wording, and a brand-new "Additions" category is appended
# as the final row.
$ws.Range("A9").Value  = "صواني"
$ws.Range("A10").Value = "عصائر فريش "
$ws.Range("A11").Value = "كوكتيلات"
$ws.Range("A12").Value = "سموذي"
$ws.Range("A13").Value = "مشروبات ساخنة"
$ws.Range("A14").Value = "قهوة"
$ws.Range("A15").Value = "مشروبات باردة"
$ws.Range("A16").Value = "فرابية"
$ws.Range("A17").Value = "ميلك شيك"
$ws.Range("A18").Value = "حلويات"
$ws.Range("A19").Value = "شيشة"

# New row 20: "Additions" category, same 14% tax rate as every other row.
$ws.Range("A20").Value = "الاضافات"
$ws.Range("C20").Value = 14

$ws.Range("B18").Select()
